$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.837.38"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").Value = "2.633.96"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "568.53"
$ws.Range("E5").Value = "  +6.45%  "
$ws.Range("D6").Value = "146.38"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.610"
$ws.Range("E8").Value = "  +4.93%  "
$ws.Range("D9").Value = "2.658.73"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("E12").Value = "  +7.01%  "
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "3.099.92"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "60.738.84"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "22.07"
$ws.Range("E16").Value = "  +6.77%  "
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "2.651.24"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").Value = "343.22"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  +3.84%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +5.18%  "
$ws.Range("D29").Value = "0.0₃0812"
$ws.Range("E29").Value = "  +11.64%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  +6.22%  "
$ws.Range("D33").Value = "160.13"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "19.23"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("E36").Value = "  +9.08%  "
$ws.Range("E37").Value = "  +5.49%  "
$ws.Range("D38").Value = "0.892"
$ws.Range("E38").Value = "  +9.33%  "
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("D40").Value = "37.42"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").Value = "300.36"
$ws.Range("E41").Value = "  +7.42%  "
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "0.0985"
$ws.Range("E44").Value = "  +4.20%  "
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0545"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "128.19"
$ws.Range("E47").Value = "  +14.68%  "
$ws.Range("D48").Value = "19.37"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "'10.70"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  +5.13%  "
$ws.Range("D51").Value = "4.65"
$ws.Range("E51").Value = "  +5.03%  "
